# Update "Forecast Comparison" sheet (sheet1) with corrected forecast output:
#  - insert a new "Week_Start_Date" column after "Week"
#  - shorten week labels (W01 -> W1, etc.)
#  - correct a few MyForecast values
#  - store is_holiday_week as a boolean
# Then refresh the dependent totals on the "Summary" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# Insert a new column B; everything from the old B onward shifts one to the right
# (ASIN: B->C, MyForecast: C->D, Amazon Mean: D->E, P70: E->F, P80: F->G,
#  P90: G->H, Product Title: H->I, is_holiday_week: I->J).
$ws.Columns.Item(2).Insert()

# New header for the inserted column.
$ws.Range("B1").Value = "Week_Start_Date"

# Corrected MyForecast values (column D after the insert) for a few weeks.
$myForecast = @{
    2  = 42
    3  = 44
    4  = 45
    5  = 50
    6  = 54
    7  = 46
    8  = 58
    9  = 58
    10 = 44
    11 = 45
    12 = 43
    13 = 51
    14 = 42
    15 = 39
    16 = 38
    17 = 39
}

# Week start dates (text, not Excel date serials) for the new column B.
$weekStart = @{
    2  = "2025-01-05"
    3  = "2025-01-12"
    4  = "2025-01-19"
    5  = "2025-01-26"
    6  = "2025-02-02"
    7  = "2025-02-09"
    8  = "2025-02-16"
    9  = "2025-02-23"
    10 = "2025-03-02"
    11 = "2025-03-09"
    12 = "2025-03-16"
    13 = "2025-03-23"
    14 = "2025-03-30"
    15 = "2025-04-06"
    16 = "2025-04-13"
    17 = "2025-04-20"
}

for ($row = 2; $row -le 17; $row++) {
    $weekNum = $row - 1

    # Shorten "W01".."W16" to "W1".."W16".
    $ws.Cells.Item($row, 1).Value = "W$weekNum"

    # Week_Start_Date as literal text.
    $cellB = $ws.Cells.Item($row, 2)
    $cellB.NumberFormat = "@"
    $cellB.Value = $weekStart[$row]

    # Corrected MyForecast value.
    $ws.Cells.Item($row, 4).Value = $myForecast[$row]

    # is_holiday_week as a real boolean instead of a 0/1 number (all weeks
    # here are non-holiday weeks).
    $cellJ = $ws.Cells.Item($row, 10)
    $cellJ.Value = [bool]([double]$cellJ.Value2)
}

# Refresh the dependent summary totals (stored as text, matching the rest
# of the column).
$summary = $wb.Worksheets.Item("Summary")
foreach ($pair in @{ B9 = "739"; B10 = "397"; B12 = "58" }.GetEnumerator()) {
    $cell = $summary.Range($pair.Key)
    $cell.NumberFormat = "@"
    $cell.Value = $pair.Value
}
